$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 13 entirely (it only held a stray "519033 - Carlos Yujiro Shigue"
# value in B/C with no label in A) - this shifts rows 14:22 up to 13:21.
$ws.Rows("13").Delete()

# After the shift, overwrite the B/C values that now hold stale/duplicated
# content with the correct values for their (new) row/label.
$ws.Range("B10").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C10").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2012" must stay literal text (matches the existing "Ativação:" cell),
# so copy it from B8/C8 instead of assigning the literal, which Excel would
# otherwise auto-convert to a date serial value.
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("C8").Copy($ws.Range("C15"))

$ws.Range("B18").Value = "519033 - Carlos Yujiro Shigue"
$ws.Range("C18").Value = "519033 - Carlos Yujiro Shigue"

$ws.Range("B19").Value = "Provas, listas de exercícios e trabalhos práticos."
$ws.Range("C19").Value = "Provas, listas de exercícios e trabalhos práticos."

$ws.Range("B20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"
$ws.Range("C20").Value = "Média ponderada de duas provas escritas, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + 2P2 + TR)/4"

$ws.Range("B21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$ws.Range("C21").Value = "Aplicação de uma prova escrita dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
